$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 845, shifting the existing rows 845:886 down to 846:887.
$ws.Rows("845").Insert()

# Populate the newly inserted row with the new record: 2026/02/20, 金, 10, 201
# Force column A to be stored as literal text (not auto-parsed as a date),
# then restore the default "Normal" style so no stray number-format style lingers.
$ws.Range("A845").NumberFormat = "@"
$ws.Range("A845").Value = "2026/02/20"
$ws.Range("A845").Style = "Normal"

$ws.Range("B845").Value = "金"
$ws.Range("C845").Value = 10
$ws.Range("D845").Value = 201
